$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Pomc -> Mc5r -> ECs)
$ws.Range("G2").Value = 1.524170333333333
$ws.Range("H2").Value = 4.572511
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5306206666666666
$ws.Range("N2").Value = 1.591862
$ws.Range("O2").Value = 0.2609258655515381
$ws.Range("P2").Value = 0.260925865551538
$ws.Range("Q2").Value = 0.8087562783868889
$ws.Range("R2").Value = 7.278806505482001
$ws.Range("S2").Value = 0.2609258655515381
$ws.Range("T2").Value = 0.260925865551538

# Row 3 (ECs -> Pomc -> Mc5r -> FAPs)
$ws.Range("G3").Value = 1.524170333333333
$ws.Range("H3").Value = 4.572511
$ws.Range("O3").Value = 0.61220989765148
$ws.Range("P3").Value = 0.6122098976514799
$ws.Range("Q3").Value = 1.897583428034778
$ws.Range("R3").Value = 17.078250852313
$ws.Range("S3").Value = 0.61220989765148
$ws.Range("T3").Value = 0.6122098976514799

# Row 4 (ECs -> Pomc -> Mc5r -> MuSCs)
$ws.Range("G4").Value = 1.524170333333333
$ws.Range("H4").Value = 4.572511
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.257992
$ws.Range("N4").Value = 0.773976
$ws.Range("O4").Value = 0.1268642367969819
$ws.Range("P4").Value = 0.1268642367969819
$ws.Range("Q4").Value = 0.3932237526373333
$ws.Range("R4").Value = 3.539013773736
$ws.Range("S4").Value = 0.1268642367969819
$ws.Range("T4").Value = 0.1268642367969819
